$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1717.3077
$ws.Range("I98").Value = 1726.1177
$ws.Range("J98").Value = 1700.6666
$ws.Range("K98").Value = 1726.1177
$ws.Range("L98").Value = 1700.6666
$ws.Range("M98").Value = -228.1177
$ws.Range("N98").Value = -4696.6666

$ws.Range("H122").Value = 1717.3077
$ws.Range("I122").Value = 1726.1177
$ws.Range("J122").Value = 1700.6666
$ws.Range("K122").Value = 5178.3531
$ws.Range("L122").Value = 5101.9998
$ws.Range("M122").Value = -2728.3531
$ws.Range("N122").Value = -10001.9998

$ws.Range("H129").Value = 1044.36
$ws.Range("I129").Value = 347.5
$ws.Range("J129").Value = 1104.9565
$ws.Range("K129").Value = 1042.5
$ws.Range("L129").Value = 3314.8695
$ws.Range("M129").Value = 3957.5
$ws.Range("N129").Value = -13314.8695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1881.1428
$ws.Range("I45").Value = 1569.6666
$ws.Range("K45").Value = 1569.6666
$ws.Range("M45").Value = -1192.6666

$ws.Range("H63").Value = 2326.6667
$ws.Range("I63").Value = 2326.6667
$ws.Range("K63").Value = 2326.6667
$ws.Range("M63").Value = -1640.6667

$ws.Range("H66").Value = 2326.6667
$ws.Range("I66").Value = 2326.6667
$ws.Range("K66").Value = 11633.3335
$ws.Range("M66").Value = -8201.333500000001

$ws.Range("H74").Value = 730.875
$ws.Range("I74").Value = 790.3158
$ws.Range("J74").Value = 677.0952
$ws.Range("K74").Value = 790.3158
$ws.Range("L74").Value = 677.0952
$ws.Range("M74").Value = 83.68420000000003
$ws.Range("N74").Value = -2425.0952

$ws.Range("H77").Value = 730.875
$ws.Range("I77").Value = 790.3158
$ws.Range("J77").Value = 677.0952
$ws.Range("K77").Value = 3951.579
$ws.Range("L77").Value = 3385.476
$ws.Range("M77").Value = 416.4210000000003
$ws.Range("N77").Value = -12121.476

$ws.Range("H122").Value = 1532.6666
$ws.Range("I122").Value = 1567.8
$ws.Range("J122").Value = 1357
$ws.Range("K122").Value = 4703.4
$ws.Range("L122").Value = 4071
$ws.Range("M122").Value = -2253.4
$ws.Range("N122").Value = -8971

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2241.6155
$ws.Range("I86").Value = 1943.4445
$ws.Range("J86").Value = 2912.5
$ws.Range("K86").Value = 1943.4445
$ws.Range("L86").Value = 2912.5
$ws.Range("M86").Value = -820.4445000000001
$ws.Range("N86").Value = -5158.5

$ws.Range("H89").Value = 2241.6155
$ws.Range("I89").Value = 1943.4445
$ws.Range("J89").Value = 2912.5
$ws.Range("K89").Value = 9717.2225
$ws.Range("L89").Value = 14562.5
$ws.Range("M89").Value = -4101.2225
$ws.Range("N89").Value = -25794.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 34057.2
$ws.Range("J125").Value = 34057.2
$ws.Range("L125").Value = 34057.2
$ws.Range("N125").Value = -38977.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 9656.909
$ws.Range("J110").Value = 14999.857
$ws.Range("L110").Value = 44999.571
$ws.Range("N110").Value = -53179.571

$ws.Range("H115").Value = 726.5
$ws.Range("I115").Value = 726.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2179.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1004.5
$ws.Range("N115").ClearContents()

$ws.Range("H118").Value = 3283.1428
$ws.Range("I118").Value = 1190
$ws.Range("J118").Value = 8516
$ws.Range("K118").Value = 3570
$ws.Range("L118").Value = 25548
$ws.Range("M118").Value = -2327
$ws.Range("N118").Value = -28034

$ws.Range("H123").Value = 4375
$ws.Range("J123").Value = 4571.4287
$ws.Range("L123").Value = 13714.2861
$ws.Range("N123").Value = -18614.2861

$ws.Range("H124").Value = 3921.4285
$ws.Range("I124").Value = 1225
$ws.Range("K124").Value = 3675
$ws.Range("M124").Value = 1235

$ws.Range("H130").Value = 2993.3333
$ws.Range("I130").Value = 1192
$ws.Range("J130").Value = 12000
$ws.Range("K130").Value = 3576
$ws.Range("L130").Value = 36000
$ws.Range("M130").Value = 1444
$ws.Range("N130").Value = -46040

$ws.Range("H139").Value = 1759.3334
$ws.Range("I139").Value = 1490
$ws.Range("J139").Value = 2500
$ws.Range("K139").Value = 4470
$ws.Range("L139").Value = 7500
$ws.Range("M139").Value = 670
$ws.Range("N139").Value = -17780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 69.666664
$ws.Range("I2").Value = 51.666668
$ws.Range("J2").Value = 87.666664
$ws.Range("K2").Value = 51.666668
$ws.Range("L2").Value = 87.666664
$ws.Range("M2").Value = 61.333332
$ws.Range("N2").Value = -313.666664

$ws.Range("H80").Value = 3464.348
$ws.Range("I80").Value = 3409
$ws.Range("J80").Value = 3833.3333
$ws.Range("K80").Value = 3409
$ws.Range("L80").Value = 3833.3333
$ws.Range("M80").Value = -2411
$ws.Range("N80").Value = -5829.3333

$ws.Range("H83").Value = 3464.348
$ws.Range("I83").Value = 3409
$ws.Range("J83").Value = 3833.3333
$ws.Range("K83").Value = 17045
$ws.Range("L83").Value = 19166.6665
$ws.Range("M83").Value = -12053
$ws.Range("N83").Value = -29150.6665

$ws.Range("H122").Value = 1510.3636
$ws.Range("I122").Value = 1604.8334
$ws.Range("J122").Value = 1397
$ws.Range("K122").Value = 4814.5002
$ws.Range("L122").Value = 4191
$ws.Range("M122").Value = -2364.5002
$ws.Range("N122").Value = -9091

$ws.Range("H132").Value = 2344.5789
$ws.Range("I132").Value = 1690.5625
$ws.Range("J132").Value = 5832.6665
$ws.Range("K132").Value = 5071.6875
$ws.Range("L132").Value = 17497.9995
$ws.Range("M132").Value = -2541.6875
$ws.Range("N132").Value = -22557.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4918.1177
$ws.Range("I122").Value = 5321.5864
$ws.Range("J122").Value = 2578
$ws.Range("K122").Value = 15964.7592
$ws.Range("L122").Value = 7734
$ws.Range("M122").Value = -13514.7592
$ws.Range("N122").Value = -12634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1943.2858
$ws.Range("I122").Value = 1760.6
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 5281.799999999999
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -2831.799999999999
$ws.Range("N122").Value = -12100

$ws.Range("H123").Value = 49320
$ws.Range("J123").Value = 49320
$ws.Range("L123").Value = 49320
$ws.Range("N123").Value = -59120
